$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "id"
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "phone"
$ws.Range("D1").Value = "email"
$ws.Range("E1").Value = "address"

$ws.Range("F3").Select()
